$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Determine the extent of the data by reading the worksheet dimension
$lastRow = $ws.UsedRange.Row + $ws.UsedRange.Rows.Count - 1
$lastCol = $ws.UsedRange.Column + $ws.UsedRange.Columns.Count - 1

for ($r = 2; $r -le $lastRow; $r++) {

    # --- Update "Förändrad" (column C) date serial from 45184 to 45186 ---
    $cCell = $ws.Cells.Item($r, 3)
    if ($cCell.Value2 -eq 45184) {
        $cCell.Value = 45186
    }

    # Friendly display name for this row's links = value of column A ("Beteckning")
    $name = $ws.Cells.Item($r, 1).Value2

    # --- Add the friendly-name second argument to every HYPERLINK() formula
    #     in columns S..Y (Artfyndslänk .. Tillsynsbegäransmaillänk) ---
    for ($col = 19; $col -le 25; $col++) {
        $cell = $ws.Cells.Item($r, $col)
        $formula = $cell.Formula
        if ($formula -ne $null -and $formula -ne "") {
            if ($formula.StartsWith("=HYPERLINK(") -and -not $formula.Contains(",")) {
                $trimmed = $formula.Substring(0, $formula.Length - 1)
                $newFormula = $trimmed + ', "' + $name + '")'
                $cell.Formula = $newFormula
            }
        }
    }
}
